# Mandarin Zones for Top Nav - update (Jan 2018 -> new IA draft)
# Approved by Marlene: remove retired nav items, repoint the last two rows
# to "learning-english" and a "+" link row to the Category Manager, and
# restyle the B/C columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the now-retired rows (old items 11-15: VOA-interviews,
#     us-china-relations, north-korea, cross-strait-relations, learning-english)
$ws.Rows("13:17").Delete()

# --- Repoint the last two remaining data rows.
# Row 11 used to be "education" -> now shows "learning-english"
$ws.Range("B11").Value = "学英语"
$ws.Range("C11").Value = "learning-english"
# Row 12 used to be "health" -> now becomes the "+" add-category link row
$ws.Range("B12").Value = "+"
$ws.Range("C12").Value = "(link to Category Manager)"

# --- Column sizing
$ws.Columns("B").ColumnWidth = 6.1667
$ws.Columns("C").ColumnWidth = 20.8

# --- C column styling (Calibri, vertically centered)
$ws.Range("C3:C5").Font.Name = "Calibri"
$ws.Range("C3:C5").VerticalAlignment = -4108

# Row 6 (economics-trade) also wraps text and is a taller row
$ws.Range("C6").Font.Name = "Calibri"
$ws.Range("C6").VerticalAlignment = -4108
$ws.Range("C6").WrapText = $true
$ws.Rows("6").RowHeight = 30

$ws.Range("C7:C10").Font.Name = "Calibri"
$ws.Range("C7:C10").VerticalAlignment = -4108

$ws.Range("C12").Font.Name = "Calibri"
$ws.Range("C12").VerticalAlignment = -4108

# --- B column styling (Microsoft JhengHei, vertically centered)
$ws.Range("B3:B10").Font.Name = "Microsoft JhengHei"
$ws.Range("B3:B10").VerticalAlignment = -4108

# The "+" link row gets a larger left-aligned font
$ws.Range("B12").Font.Name = "Microsoft JhengHei"
$ws.Range("B12").Font.Size = 16
$ws.Range("B12").VerticalAlignment = -4108
$ws.Range("B12").HorizontalAlignment = -4131
$ws.Rows("12").RowHeight = 17

# --- Reset selection/view back to the top of the sheet
$ws.Range("A1").Select()
